$d = $word.ActiveDocument

# Locate the paragraph "LOB1019: Física II (Requisito fraco)" - it is kept as-is.
# Everything from the paragraph immediately following it, through the end of the
# copyright/footer paragraph ("© 2020 ... Powered by Jekyll ..."), must be removed
# (the blank paragraph, the "Ver no Jupiter ..." paragraph, and the copyright paragraph).

$wdParagraph = 4

$anchor = $d.Content.Duplicate
[void]$anchor.Find.Execute("LOB1019: Física II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$anchor.Expand($wdParagraph)

$footer = $d.Content.Duplicate
[void]$footer.Find.Execute("Powered by Jekyll and Github pages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$footer.Expand($wdParagraph)

$toDelete = $d.Range($anchor.End, $footer.End)
[void]$toDelete.Delete()
